# Auto-generated Excel COM-interop script replicating the diff
$wb = $excel.ActiveWorkbook

# --- INDI sheet: new rows 106-121 (2022 Q1 .. 2025 Q4) ---
$wsIndi = $wb.Worksheets.Item("INDI")
$indiDates = @("2022-01-01","2022-04-01","2022-07-01","2022-10-01","2023-01-01","2023-04-01","2023-07-01","2023-10-01","2024-01-01","2024-04-01","2024-07-01","2024-10-01","2025-01-01","2025-04-01","2025-07-01","2025-10-01")
$indiB     = @("652337.430519783","1391602.35233014","1377764.88603498","1497787.69083112","674382.277702118","1438670.81950154","1424390.5208906","1548491.52084353","697216.302293026","1487388.90415494","1472628.57410196","1600934.6945855","720829.756119316","1537764.91054443","1522505.17145261","1655157.23856895")
$indiC     = @("65.2337430519783","139.160235233014","137.776488603498","149.778769083112","67.4382277702118","143.867081950154","142.43905208906","154.849152084353","69.7216302293026","148.738890415494","147.262857410196","160.09346945855","72.0829756119316","153.776491054443","152.250517145261","165.515723856895")
for ($i = 0; $i -lt $indiDates.Length; $i++) {
    $r = 106 + $i
    $cellA = $wsIndi.Range("A$r")
    $cellA.NumberFormat = "@"
    $cellA.Value = $indiDates[$i]
    $cellA.Style = "Normal"
    $wsIndi.Range("B$r").Value = [double]$indiB[$i]
    $wsIndi.Range("C$r").Value = [double]$indiC[$i]
}

# --- ETALONNAGE sheet: new rows 28-31 (annual 2022..2025) ---
$wsEtal = $wb.Worksheets.Item("ETALONNAGE")
$etalDates = @("2022-01-01","2023-01-01","2024-01-01","2025-01-01")
$etalB     = @("1488686.75033894","1526741.06948686","1566119.30129854","1606836.24024701")
$etalC     = @("122.987308992901","127.148378473445","131.454211878386","135.906426917133")
$etalD     = @("2.5185612948861","2.55623415330679","2.57923446212907","2.59986189524031")
$etalE     = @("3.36085799827193","3.38333240609754","3.38646348198632","3.38689417031837")
$etalF     = @("Acceptable","Acceptable","Acceptable","Acceptable")
for ($i = 0; $i -lt $etalDates.Length; $i++) {
    $r = 28 + $i
    $cellA = $wsEtal.Range("A$r")
    $cellA.NumberFormat = "@"
    $cellA.Value = $etalDates[$i]
    $cellA.Style = "Normal"
    $wsEtal.Range("B$r").Value = [double]$etalB[$i]
    $wsEtal.Range("C$r").Value = [double]$etalC[$i]
    $wsEtal.Range("D$r").Value = [double]$etalD[$i]
    $wsEtal.Range("E$r").Value = [double]$etalE[$i]
    $wsEtal.Range("F$r").Value = $etalF[$i]
}

# --- PREVISION sheet: new rows 102-117 (2022 Q1 .. 2025 Q4) ---
$wsPrev = $wb.Worksheets.Item("PREVISION")
$prevDates = @("2022-01-01","2022-04-01","2022-07-01","2022-10-01","2023-01-01","2023-04-01","2023-07-01","2023-10-01","2024-01-01","2024-04-01","2024-07-01","2024-10-01","2025-01-01","2025-04-01","2025-07-01","2025-10-01")
$prevB     = @("240127.952162735","409148.341037331","405984.639006779","433425.818132099","245168.132663001","419909.74792995","416644.799799372","445018.389094539","250388.745195888","431048.311887751","427673.612707702","457008.631507203","255787.561006894","442565.931537618","439077.050949171","469405.696753326")
$prevC     = @("16.3084357629946","34.7900588082535","34.4441221508745","37.444692270778","16.8595569425529","35.9667704875385","35.609763022265","38.7122880210883","17.4304075573257","37.1847226038735","36.815714352549","40.0233673646375","18.0207439029829","38.4441227636107","38.0626292863152","41.3789309642237")
$prevD     = @("12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233","12181.1365322233")
for ($i = 0; $i -lt $prevDates.Length; $i++) {
    $r = 102 + $i
    $cellA = $wsPrev.Range("A$r")
    $cellA.NumberFormat = "@"
    $cellA.Value = $prevDates[$i]
    $cellA.Style = "Normal"
    $wsPrev.Range("B$r").Value = [double]$prevB[$i]
    $wsPrev.Range("C$r").Value = [double]$prevC[$i]
    $wsPrev.Range("D$r").Value = [double]$prevD[$i]
}

# --- VATRIM sheet: new rows 102-117 (2022 Q1 .. 2025 Q4) ---
$wsVat = $wb.Worksheets.Item("VATRIM")
$vatDates = @("2022-01-01","2022-04-01","2022-07-01","2022-10-01","2023-01-01","2023-04-01","2023-07-01","2023-10-01","2024-01-01","2024-04-01","2024-07-01","2024-10-01","2025-01-01","2025-04-01","2025-07-01","2025-10-01")
$vatB     = @("240127.952162735","409148.341037331","405984.639006779","433425.818132099","245168.132663001","419909.74792995","416644.799799372","445018.389094539","250388.745195888","431048.311887751","427673.612707702","457008.631507203","255787.561006894","442565.931537618","439077.050949171","469405.696753326")
for ($i = 0; $i -lt $vatDates.Length; $i++) {
    $r = 102 + $i
    $cellA = $wsVat.Range("A$r")
    $cellA.NumberFormat = "@"
    $cellA.Value = $vatDates[$i]
    $cellA.Style = "Normal"
    $wsVat.Range("B$r").Value = [double]$vatB[$i]
}
